# Applies the "Eigen_Terrasse" sheet rework:
#   - removes the leading "Zeile" (row-number) column, shifting A:F -> A:E
#   - inserts a new row for the "Oberfläche" / P_OF surface-treatment factor
#     between the "Anzahl Saeulen" row and the "Dach-Eindeckung" row
#   - updates the internal "Traegerpreis/m" formula constants (60/110 -> 90/160)
#   - updates the final price formula to multiply by the new P_OF factor

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eigen_Terrasse")

# Drop column A ("Zeile") entirely; B:F shift left to become A:E.
$ws.Columns.Item(1).Delete()

# Make room for the new "Oberfläche" row above the "Dach-Eindeckung" row
# (old row 6, now at row 6 after the column delete).
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "Auswahl"
$ws.Range("B6").Value = "Oberfläche"
$ws.Range("C6").Value = "P_OF"
$ws.Range("D6").Value = "Roh keine Oberflächenbehandlung:1, feuerverzinkt:1.3, KTL und Pulverbeschichtet:1.6"

# Update the internal helper formula text for the beam price per metre.
$ws.Range("E11").Value = "(L <= 5) * 90 + (L > 5) * 160"

# Update the overall price formula to fold in the new P_OF surface factor.
$ws.Range("E12").Value = "( (L * P_Trager * P_OF) + (N_Col * H * 90* P_OF) + (N_Spar * B * 110* P_OF) + (L * B * P_Dach) + (L * P_Wand) ) * (1 - (Rabatt / 100))"

$ws.Range("E12").Select()
